$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 6192
    $ws.Range("F6").Value = 22
    $ws.Range("F8").Value = 1862
    $ws.Range("F9").Value = 1397
    $ws.Range("F10").Value = 292
    $ws.Range("F11").Value = 949
    $ws.Range("F12").Value = 206
    $ws.Range("F13").Value = 5566
}
